# Daily attendance processing - 2026-02-19 10:52:05 UTC
# Reorders the comma-separated tokens in column G for specific rows
# (swapping cohort years / reviewer emails to the order produced by the
# latest processing run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Session Analysis Results")

$updates = @{
    15  = "2025/2026, 2023/2024"
    19  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    21  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    22  = "2025/2026, Eman_mohamed@med.asu.edu.eg"
    37  = "2025/2026, 2023/2024"
    41  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    43  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    44  = "2025/2026, Eman_mohamed@med.asu.edu.eg"
    60  = "2026/2027, 2025/2026"
    63  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    64  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    65  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    66  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    82  = "2026/2027, 2025/2026"
    85  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    86  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    87  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    88  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    89  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    104 = "2026/2027, 2025/2026"
    106 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    107 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    110 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    111 = "marina_atef@med.asu.edu.eg, 2025/2026, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
    126 = "2026/2027, 2025/2026"
    128 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    129 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    132 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    133 = "marina_atef@med.asu.edu.eg, 2025/2026, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
    150 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    153 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    155 = "marina_atef@med.asu.edu.eg, 2025/2026, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
    172 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    175 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    177 = "marina_atef@med.asu.edu.eg, 2025/2026, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
